# Reorder the "open/closed world" and "distance from concept" related
# columns (Q:V) on the LanguageCandidates sheet.
#
# Old column order (Q..V):
#   Q = DimensionalityWhileEditing
#   R = IsOpenWorld
#   S = IsClosedWorld
#   T = IsOpenClosedWorldConflicted
#   U = DistanceFromConcept
#   V = IsDescriptionOf
#
# New column order (Q..V):
#   Q = IsOpenWorld
#   R = IsClosedWorld
#   S = IsDescriptionOf
#   T = DistanceFromConcept
#   U = IsOpenClosedWorldConflicted
#   V = DimensionalityWhileEditing
#
# This is achieved with three column cut/insert moves (Excel automatically
# re-writes the relative formula references as columns shift):
#   1. Move DimensionalityWhileEditing (Q) to the end of the block (after V).
#   2. Swap IsOpenClosedWorldConflicted and DistanceFromConcept (now at S/T).
#   3. Move IsDescriptionOf (now at U) to right after IsClosedWorld (R).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LanguageCandidates")

# Step 1: move Q (DimensionalityWhileEditing) to the end of the Q:V block.
$ws.Columns("Q:Q").Cut() | Out-Null
$ws.Columns("W:W").Insert() | Out-Null
# Now: Q=IsOpenWorld, R=IsClosedWorld, S=IsOpenClosedWorldConflicted,
#      T=DistanceFromConcept, U=IsDescriptionOf, V=DimensionalityWhileEditing

# Step 2: swap IsOpenClosedWorldConflicted (S) and DistanceFromConcept (T).
$ws.Columns("T:T").Cut() | Out-Null
$ws.Columns("S:S").Insert() | Out-Null
# Now: Q=IsOpenWorld, R=IsClosedWorld, S=DistanceFromConcept,
#      T=IsOpenClosedWorldConflicted, U=IsDescriptionOf, V=DimensionalityWhileEditing

# Step 3: move IsDescriptionOf (U) to right after IsClosedWorld (R).
$ws.Columns("U:U").Cut() | Out-Null
$ws.Columns("S:S").Insert() | Out-Null
# Now: Q=IsOpenWorld, R=IsClosedWorld, S=IsDescriptionOf,
#      T=DistanceFromConcept, U=IsOpenClosedWorldConflicted, V=DimensionalityWhileEditing
